$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

$ws.Cells.Item($row, 1).Value = "46SS1B"
$ws.Cells.Item($row, 2).Value = "Cinta Flex de sensor y cabezal Epson"
$ws.Cells.Item($row, 3).Value = "L1110 L1118 L1119 L1210 L1250  L3100 L3101 L3106 L3108 L3109 L3110 L3115 L3116 L3117 L3118 L3119 L3150 L3151 L3153 L3156 L3158 L3160 L3161 L3163 L3166  L3167 L3168 L3169 L3210 L3216 L3250 L3256 L5190 L5290"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 150000
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 8
$ws.Cells.Item($row, 8).Formula = "=(E25-D25)*G25"
$ws.Cells.Item($row, 9).Formula = "=D25*F25"
$ws.Cells.Item($row, 10).Value = 0
